{"js": "// skillkaart kr1 stelt de opdracht vast: foutjes verbetert\n//\n// 1. \"mijn vriend\" -> \"een vriend\" (appears twice: intro paragraph + \"De\n//    website wordt gemaakt voor ...\" paragraph).\n// 2. \"Alle bedrijven moeten op de website te komen staan (...)\" ->\n//    \"Alle nep bedrijven moeten op de website komen te staan (...)\".\n// 3. \"Het moet aan de logo kleuren combineren.\" ->\n//    \"De website moet aan de logo kleuren combineren.\".\n// 4. The \"_GoBack\" bookmark (Word's \"last edit location\" marker) follows\n//    the edit and ends up collapsed right after the \"een\" that replaced\n//    \"mijn\" in the \"De website wordt gemaakt ...\" paragraph, instead of\n//    sitting in the trailing empty paragraph at the end of the document.\n\n// --- 1. \"mijn vriend\" -> \"een vriend\" (both occurrences) -----------------\nconst mijnVriend = context.document.body.search(\"mijn vriend\", { matchCase: false });\nmijnVriend.load(\"text\");\nawait context.sync();\n\nfor (const r of mijnVriend.items) {\n  r.insertText(\"een vriend\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. \"Alle bedrijven moeten op de website te komen staan\" -------------\nconst bedrijvenZin = context.document.body.search(\n  \"Alle bedrijven moeten op de website te komen staan (hoofdkantoor, en 4 winkels).\",\n  { matchCase: false }\n);\nbedrijvenZin.load(\"text\");\nawait context.sync();\n\nbedrijvenZin.items[0].insertText(\n  \"Alle nep bedrijven moeten op de website komen te staan (hoofdkantoor, en 4 winkels).\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 3. \"Het moet aan de logo kleuren combineren.\" ------------------------\nconst logoZin = context.document.body.search(\"Het moet aan de logo kleuren combineren.\", {\n  matchCase: false,\n});\nlogoZin.load(\"text\");\nawait context.sync();\n\nlogoZin.items[0].insertText(\n  \"De website moet aan de logo kleuren combineren.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 4. Move the \"_GoBack\" bookmark --------------------------------------\n// Word always keeps a single \"_GoBack\" bookmark that marks the location of\n// the most recent edit. Originally it sat in the empty paragraph at the\n// end of the document; after this edit it should collapse to the point\n// right after \"een\" in \"De website wordt gemaakt voor een vriend(en) ...\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst doelParagraaf = context.document.body.search(\"De website wordt gemaakt voor een\", {\n  matchCase: false,\n});\ndoelParagraaf.load(\"text\");\nawait context.sync();\n\nconst paragraaf = doelParagraaf.items[0].paragraphs.getFirst();\n\n// \"een\" occurs as a whole word three times in that paragraph; the first\n// one is the word that used to be \"mijn\".\nconst eenWoord = paragraaf.search(\"een\", { matchCase: true, matchWholeWord: true });\neenWoord.load(\"text\");\nawait context.sync();\n\nconst eenEinde = eenWoord.items[0].getRange(Word.RangeLocation.end);\neenEinde.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# skillkaart kr1 stelt de opdracht vast: foutjes verbetert\n#\n# 1. \"mijn vriend\" -> \"een vriend\" (appears twice: intro paragraph + \"De\n#    website wordt gemaakt voor ...\" paragraph).\n# 2. \"Alle bedrijven moeten op de website te komen staan (...)\" ->\n#    \"Alle nep bedrijven moeten op de website komen te staan (...)\".\n# 3. \"Het moet aan de logo kleuren combineren.\" ->\n#    \"De website moet aan de logo kleuren combineren.\".\n# 4. The \"_GoBack\" bookmark (Word's \"last edit location\" marker) follows\n#    the edit and ends up collapsed right after the \"een\" that replaced\n#    \"mijn\" in the \"De website wordt gemaakt ...\" paragraph, instead of\n#    sitting in the trailing empty paragraph at the end of the document.\n\n$d = $word.ActiveDocument\n\n# --- 1a. \"mijn vriend\" -> \"een vriend\" in the intro paragraph -------------\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"mijn vriend te helpen\"\n$find.Forward = $true\n$find.Wrap = 0\nif ($find.Execute()) {\n    $range.Text = \"een vriend te helpen\"\n}\n\n# --- 1b. \"mijn vriend(en)\" -> \"een vriend(en)\" in \"De website wordt ------\n#         gemaakt voor ...\" -- only the \"mijn\" word is replaced, and the\n#         \"_GoBack\" bookmark is re-created right after it.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"mijn vriend(en)\"\n$find2.Forward = $true\n$find2.Wrap = 0\n$find2.Execute()\n\n$mijnStart = $range2.Start\n$mijnRange = $d.Range($mijnStart, $mijnStart + 4)\n$mijnRange.Text = \"een\"\n\n$goBackPoint = $d.Range($mijnStart + 3, $mijnStart + 3)\n$d.Bookmarks.Add(\"_GoBack\", $goBackPoint)\n\n# --- 2. \"Alle bedrijven moeten op de website te komen staan\" --------------\n$range3 = $d.Content\n$find3 = $range3.Find\n$find3.ClearFormatting()\n$find3.Text = \"Alle bedrijven moeten op de website te komen staan\"\n$find3.Forward = $true\n$find3.Wrap = 0\nif ($find3.Execute()) {\n    $range3.Text = \"Alle nep bedrijven moeten op de website komen te staan\"\n}\n\n# --- 3. \"Het moet aan de logo kleuren combineren.\" -------------------------\n$range4 = $d.Content\n$find4 = $range4.Find\n$find4.ClearFormatting()\n$find4.Text = \"Het moet aan de logo kleuren combineren.\"\n$find4.Forward = $true\n$find4.Wrap = 0\nif ($find4.Execute()) {\n    $range4.Text = \"De website moet aan de logo kleuren combineren.\"\n}\n"}
